$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update alpha_distance_range row (row 2)
$ws.Range("B2").Value = 5.6
$ws.Range("C2").Value = 9.1999999999999993

# Update beta_distance_range row (row 3)
$ws.Range("B3").Value = 5.7
$ws.Range("C3").Value = 9

# Update ratio_threshold_range row (row 4)
$ws.Range("B4").Value = 0.9
$ws.Range("C4").Value = 13

# Remove the theta_threshold_range row (row 5) entirely; this shifts
# the pie_threshold_range row up from row 6 to row 5.
$ws.Rows(5).Delete()

# Set the new values for the (now-shifted) pie_threshold_range row (row 5)
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 15

# Move the selection like the saved file shows (selection moved past the data)
$ws.Range("D10").Select()
